# Insert a new weekly price-report row for Coliflor (Agrícola del Norte S.A. de
# Arica) right before the current row 46. This pushes all subsequent rows
# (old 46..133) down by one (to 47..134), matching the target diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 46 - Excel shifts rows 46..133 down to 47..134.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new record's data.
$ws.Range("A46").Value = 1
$ws.Range("B46").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C46").Value = "Arica y Parinacota"
$ws.Range("D46").Value = 44894
$ws.Range("E46").Value = 15
$ws.Range("F46").Value = 100112008
$ws.Range("G46").Value = "Coliflor"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Tercera"
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 350
$ws.Range("L46").Value = 450
$ws.Range("M46").Value = 400
$ws.Range("N46").Value = "$/unidad"
$ws.Range("O46").Value = "Región de Arica y Parinacota"
$ws.Range("P46").Value = 400
$ws.Range("Q46").Value = 1
$ws.Range("R46").Value = "Hortaliza"
